$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 250.60294
$ws.Range("I33").Value = 197.56923
$ws.Range("K33").Value = 197.56923
$ws.Range("M33").Value = 31.43077

$ws.Range("H112").Value = 2128.889
$ws.Range("J112").Value = 2161.1428
$ws.Range("L112").Value = 6483.428400000001
$ws.Range("N112").Value = -8699.4284

$ws.Range("H115").Value = 1931.9166
$ws.Range("I115").Value = 1753
$ws.Range("J115").Value = 3900
$ws.Range("K115").Value = 5259
$ws.Range("L115").Value = 11700
$ws.Range("M115").Value = -3692
$ws.Range("N115").Value = -14834

$ws.Range("H118").Value = 1478.7778
$ws.Range("J118").Value = 2024.5
$ws.Range("L118").Value = 6073.5
$ws.Range("N118").Value = -9387.5

$ws.Range("H123").Value = 35720
$ws.Range("J123").Value = 35720
$ws.Range("L123").Value = 35720
$ws.Range("N123").Value = -45520

$ws.Range("H124").Value = 54332.668
$ws.Range("J124").Value = 54332.668
$ws.Range("L124").Value = 54332.668
$ws.Range("N124").Value = -64152.668

$ws.Range("H129").Value = 1017.8333
$ws.Range("J129").Value = 965.9286
$ws.Range("L129").Value = 2897.7858
$ws.Range("N129").Value = -12897.7858

$ws.Range("H132").Value = 15146.272
$ws.Range("I132").Value = 1949.459
$ws.Range("J132").Value = 176147.4
$ws.Range("K132").Value = 5848.377
$ws.Range("L132").Value = 528442.2
$ws.Range("M132").Value = -3318.377
$ws.Range("N132").Value = -533502.2

$ws.Range("H138").Value = 2166.457
$ws.Range("I138").Value = 1835.381
$ws.Range("J138").Value = 2308.347
$ws.Range("K138").Value = 5506.143
$ws.Range("L138").Value = 6925.041000000001
$ws.Range("M138").Value = -366.143
$ws.Range("N138").Value = -17205.041

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 50000
$ws.Range("J7").Value = 50000
$ws.Range("L7").Value = 50000
$ws.Range("N7").Value = -50228

$ws.Range("H61").Value = 2249.7273
$ws.Range("I61").Value = 1333.2667
$ws.Range("J61").Value = 4213.5713
$ws.Range("K61").Value = 1333.2667
$ws.Range("L61").Value = 4213.5713
$ws.Range("M61").Value = -1121.2667
$ws.Range("N61").Value = -4637.5713

$ws.Range("H74").Value = 2096.5312
$ws.Range("I74").Value = 1803.6
$ws.Range("K74").Value = 1803.6
$ws.Range("M74").Value = -929.5999999999999

$ws.Range("H77").Value = 2096.5312
$ws.Range("I77").Value = 1803.6
$ws.Range("K77").Value = 9018
$ws.Range("M77").Value = -4650

$ws.Range("H124").Value = 28000
$ws.Range("J124").Value = 28000
$ws.Range("L124").Value = 28000
$ws.Range("N124").Value = -37820

$ws.Range("H136").Value = 2249.7273
$ws.Range("I136").Value = 1333.2667
$ws.Range("J136").Value = 4213.5713
$ws.Range("K136").Value = 3999.800099999999
$ws.Range("L136").Value = 12640.7139
$ws.Range("M136").Value = -1449.800099999999
$ws.Range("N136").Value = -17740.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 10032.5
$ws.Range("J46").Value = 10032.5
$ws.Range("L46").Value = 10032.5
$ws.Range("N46").Value = -10628.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 195189.53
$ws.Range("I4").Value = 3999
$ws.Range("J4").Value = 204749.05
$ws.Range("K4").Value = 3999
$ws.Range("L4").Value = 204749.05
$ws.Range("M4").Value = -3887
$ws.Range("N4").Value = -204973.05

$ws.Range("H118").Value = 48742
$ws.Range("J118").Value = 48742
$ws.Range("L118").Value = 48742
$ws.Range("N118").Value = -52056

$ws.Range("H124").Value = 35000
$ws.Range("J124").Value = 35000
$ws.Range("L124").Value = 35000
$ws.Range("N124").Value = -39910

$ws.Range("H140").Value = 47499.5
$ws.Range("J140").Value = 47499.5
$ws.Range("L140").Value = 47499.5
$ws.Range("N140").Value = -57859.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 11606693
$ws.Range("I33").Value = 57.375
$ws.Range("J33").Value = 17796900
$ws.Range("K33").Value = 344.25
$ws.Range("L33").Value = 106781400
$ws.Range("M33").Value = -61.25
$ws.Range("N33").Value = -106781966

$ws.Range("H80").Value = 77099140
$ws.Range("I80").Value = 1002500
$ws.Range("J80").Value = 90934890
$ws.Range("K80").Value = 3007500
$ws.Range("L80").Value = 272804670
$ws.Range("M80").Value = -3006564
$ws.Range("N80").Value = -272806542

$ws.Range("H83").Value = 77099140
$ws.Range("I83").Value = 1002500
$ws.Range("J83").Value = 90934890
$ws.Range("K83").Value = 9022500
$ws.Range("L83").Value = 818414010
$ws.Range("M83").Value = -9017820
$ws.Range("N83").Value = -818423370

$ws.Range("H136").Value = 35716756
$ws.Range("I136").Value = 71430510
$ws.Range("J136").Value = 2999.8572
$ws.Range("K136").Value = 214291530
$ws.Range("L136").Value = 8999.5716
$ws.Range("M136").Value = -214286430
$ws.Range("N136").Value = -19199.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1104.5
$ws.Range("I102").Value = 1119.4286
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 1119.4286
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = 502.5714
$ws.Range("N102").Value = -4244

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1975.8235
$ws.Range("I16").Value = 1755.9062
$ws.Range("J16").Value = 5494.5
$ws.Range("K16").Value = 1755.9062
$ws.Range("L16").Value = 5494.5
$ws.Range("M16").Value = -1585.9062
$ws.Range("N16").Value = -5834.5

$ws.Range("H122").Value = 202680.8
$ws.Range("I122").Value = 202680.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 608042.3999999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -605592.3999999999
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 157144.4
$ws.Range("I136").Value = 218124.89
$ws.Range("K136").Value = 654374.67
$ws.Range("M136").Value = -651824.67

$ws.Range("H139").Value = 45582.715
$ws.Range("I139").Value = 60958
$ws.Range("J139").Value = 44400
$ws.Range("K139").Value = 60958
$ws.Range("L139").Value = 44400
$ws.Range("M139").Value = -55818
$ws.Range("N139").Value = -54680

$ws.Range("H141").Value = 14416.556
$ws.Range("J141").Value = 14416.556
$ws.Range("L141").Value = 14416.556
$ws.Range("N141").Value = -24776.556
